$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final roster data (Player, Position, Team) for rows 2-19
$data = @(
    @("Paul George", "SG,SF,PF", "Philadelphia 76ers"),
    @("Giannis Antetokounmpo", "PF,C", "Milwaukee Bucks"),
    @("Bobby Portis", "PF,C", "Milwaukee Bucks"),
    @("Clint Capela", "C", "Atlanta Hawks"),
    @("Tyus Jones", "PG", "Phoenix Suns"),
    @("Jaren Jackson Jr.", "PF,C", "Memphis Grizzlies"),
    @("Andrew Wiggins", "SF,PF", "Golden State Warriors"),
    @("Aaron Gordon", "PF,C", "Denver Nuggets"),
    @("Jayson Tatum", "SF,PF", "Boston Celtics"),
    @("Anfernee Simons", "PG,SG", "Portland Trail Blazers"),
    @("Fred VanVleet", "PG", "Houston Rockets"),
    @("James Harden", "PG,SG", "LA Clippers"),
    @("Anthony Edwards", "SG,SF", "Minnesota Timberwolves"),
    @("Nicolas Claxton", "C", "Brooklyn Nets"),
    @("Ivica Zubac", "C", "LA Clippers"),
    @("Zion Williamson", "PF,C", "New Orleans Pelicans"),
    @("Kyle Kuzma", "PF", "Washington Wizards"),
    @("Jaden Ivey", "PG,SG", "Detroit Pistons")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
